$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Regular rows: only the Price (D) and/or Volume(1h) (E) cells change.
# Values that would otherwise be re-interpreted by Excel's type inference
# as numbers (and so lose significant trailing zeros / exact text form)
# are entered with a leading apostrophe to force them to stay text, just
# like the source data which is plain text in the workbook.
# -------------------------------------------------------------------------

# Row 2
$ws.Range("D2").Value = "42.927.89"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").Value = "2.208.11"
$ws.Range("E3").Value = "  -1.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'254.72"
$ws.Range("E5").Value = "  +3.80%  "

# Row 6
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("D7").Value = "'76.06"
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  -4.74%  "

# Row 10
$ws.Range("D10").Value = "'41.86"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11
$ws.Range("E11").Value = "  -2.77%  "

# Row 12
$ws.Range("E12").Value = "  -1.29%  "

# Row 13
$ws.Range("E13").Value = "  +0.59%  "

# Row 14
$ws.Range("D14").Value = "2.537.42"
$ws.Range("E14").Value = "  -1.77%  "

# Row 15
$ws.Range("D15").Value = "'14.40"
$ws.Range("E15").Value = "  -1.73%  "

# Row 16
$ws.Range("D16").Value = "2.205.94"
$ws.Range("E16").Value = "  -2.28%  "

# Row 17
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  -3.62%  "

# Row 18
$ws.Range("D18").Value = "42.816.39"

# Row 19
$ws.Range("E19").Value = "  -2.69%  "

# Row 20
$ws.Range("D20").Value = "'71.22"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("E21").Value = "  -1.22%  "

# Row 23
$ws.Range("E23").Value = "  -0.81%  "

# Row 24
$ws.Range("D24").Value = "'9.25"
$ws.Range("E24").Value = "  -8.95%  "

# Row 25
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("E26").Value = "  -2.30%  "

# Row 27
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").Value = "'39.13"
$ws.Range("E28").Value = "  +0.91%  "

# Row 29
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +3.69%  "

# Row 30
$ws.Range("E30").Value = "  -2.81%  "

# Row 31
$ws.Range("D31").Value = "'173.56"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").Value = "'20.19"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").Value = "'0.0853"
$ws.Range("E33").Value = "  +6.86%  "

# Row 34
$ws.Range("D34").Value = "'5.20"
$ws.Range("E34").Value = "  -2.77%  "

# Row 35
$ws.Range("E35").Value = "  -1.30%  "

# Row 38
$ws.Range("D38").Value = "'4.28"
$ws.Range("E38").Value = "  -1.29%  "

# Row 39
$ws.Range("D39").Value = "'12.40"
$ws.Range("E39").Value = "  -4.43%  "

# Row 40
$ws.Range("D40").Value = "'2.10"
$ws.Range("E40").Value = "  -1.99%  "

# Row 41
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  +16.09%  "

# Row 42
$ws.Range("D42").Value = "'0.197"
$ws.Range("E42").Value = "  -3.19%  "

# Row 43
$ws.Range("E43").Value = "  -5.39%  "

# Row 44
$ws.Range("D44").Value = "'59.70"
$ws.Range("E44").Value = "  -0.42%  "

# Row 45
$ws.Range("D45").Value = "'101.50"
$ws.Range("E45").Value = "  -3.64%  "

# Row 46
$ws.Range("D46").Value = "'0.0976"
$ws.Range("E46").Value = "  -1.86%  "

# Row 47
$ws.Range("D47").Value = "'8.27"
$ws.Range("E47").Value = "  -4.82%  "

# Row 48
$ws.Range("D48").Value = "'0.451"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49
$ws.Range("E49").Value = "  -0.66%  "

# Row 50
$ws.Range("E50").Value = "  -1.67%  "

# -------------------------------------------------------------------------
# Row 36 and 37: Kaspa/VeChain swap places (along with new price/volume).
# -------------------------------------------------------------------------
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0354"
$ws.Range("E36").Value = "  +6.25%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.107"
$ws.Range("E37").Value = "  -2.46%  "

# -------------------------------------------------------------------------
# Row 51: RocketPoolETH replaced by HuobiToken.
# -------------------------------------------------------------------------
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.65"
$ws.Range("E51").Value = "  -0.63%  "
